$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Change C19, C24, C25 from numeric 101 to text "141m"
$ws.Range("C19").Value = "141m"
$ws.Range("C24").Value = "141m"
$ws.Range("C25").Value = "141m"

# Update the selected cell/active cell as recorded in the sheet view
$ws.Range("D30").Select()
